# CPL_JLCPCB_teensy_arena_12-12_v0p1_r1.xlsx -- "Fixed orientations in 12-12 teensy 4.1 arena fab files."
# Column E holds part "Rotation" (degrees). This pick-and-place rotation data was
# corrected for a batch of parts (row 110, and rows 345-426).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 110: single isolated rotation fix ---
$ws.Range("E110").Value = 0

# --- Rows 345-426: bulk rotation corrections ---
$ws.Range("E345").Value = 180
$ws.Range("E346").Value = 90
$ws.Range("E347").Value = 0
$ws.Range("E348").Value = 90
$ws.Range("E349").Value = 90
$ws.Range("E350").Value = 0
$ws.Range("E351").Value = 0
$ws.Range("E352").Value = 0
$ws.Range("E353").Value = 90
$ws.Range("E354").Value = 90
$ws.Range("E355").Value = 90
$ws.Range("E356").Value = 90
$ws.Range("E357").Value = 0
$ws.Range("E358").Value = 0
$ws.Range("E359").Value = 270
$ws.Range("E360").Value = 0
$ws.Range("E361").Value = 0
$ws.Range("E362").Value = 0
$ws.Range("E363").Value = 0
$ws.Range("E364").Value = 0
$ws.Range("E365").Value = 0
$ws.Range("E366").Value = 0
$ws.Range("E367").Value = 0
$ws.Range("E368").Value = 90
$ws.Range("E369").Value = 90
$ws.Range("E370").Value = 90
$ws.Range("E371").Value = 90
$ws.Range("E372").Value = 90
$ws.Range("E373").Value = 90
$ws.Range("E374").Value = 90
$ws.Range("E375").Value = 90
$ws.Range("E376").Value = 0
$ws.Range("E377").Value = 0
$ws.Range("E378").Value = 180
$ws.Range("E379").Value = 90
$ws.Range("E380").Value = 90
$ws.Range("E381").Value = 90
$ws.Range("E382").Value = 90
$ws.Range("E383").Value = -240
$ws.Range("E384").Value = -240
$ws.Range("E385").Value = -240
$ws.Range("E386").Value = -240
$ws.Range("E387").Value = -210
$ws.Range("E388").Value = -210
$ws.Range("E389").Value = -210
$ws.Range("E390").Value = -210
$ws.Range("E391").Value = 180
$ws.Range("E392").Value = 180
$ws.Range("E393").Value = 180
$ws.Range("E394").Value = 180
$ws.Range("E395").Value = -150
$ws.Range("E396").Value = -150
$ws.Range("E397").Value = -150
$ws.Range("E398").Value = -150
$ws.Range("E399").Value = -120
$ws.Range("E400").Value = -120
$ws.Range("E401").Value = -120
$ws.Range("E402").Value = -120
$ws.Range("E403").Value = -90
$ws.Range("E404").Value = -90
$ws.Range("E405").Value = -90
$ws.Range("E406").Value = -90
$ws.Range("E407").Value = -60
$ws.Range("E408").Value = -60
$ws.Range("E409").Value = -60
$ws.Range("E410").Value = -60
$ws.Range("E411").Value = -30
$ws.Range("E412").Value = -30
$ws.Range("E413").Value = -30
$ws.Range("E414").Value = -30
$ws.Range("E415").Value = 0
$ws.Range("E416").Value = 0
$ws.Range("E417").Value = 0
$ws.Range("E418").Value = 0
$ws.Range("E419").Value = 30
$ws.Range("E420").Value = 30
$ws.Range("E421").Value = 30
$ws.Range("E422").Value = 30
$ws.Range("E423").Value = 60
$ws.Range("E424").Value = 60
$ws.Range("E425").Value = 60
$ws.Range("E426").Value = 60

# The corrected value for E345 was (re-)pasted in from the source data together with
# its neighbouring (empty) cell, which is why F345 picks up E-column formatting here.
$ws.Range("E345").Copy($ws.Range("F345"))
$ws.Range("F345").Value = ""

# --- Restore the view/selection state left behind by the edit ---
$ws.Range("G90").Select()
